# Fix Training Data Issue (#48)
#
# NBA box-score data was taken from 1 day off due to how the stats were
# shown, so the "Date" label column was storing "5-22-2007-08" (a mangled
# combination of the day-of-month and season label) instead of the actual
# game date. Correct it to the proper ISO date string "2008-05-22" for
# every data row (rows 2-31, column BF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = 58          # column BF
$firstRow = 2
$lastRow = 31
$correctDate = "2008-05-22"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)

    # Force the cell to stay plain text: without this, Excel would
    # auto-recognize the ISO-formatted "2008-05-22" string as a date and
    # silently convert it into a date serial number instead of keeping it
    # as the literal text the original "5-22-2007-08" value was.
    $cell.NumberFormat = "@"
    $cell.Value = $correctDate
}
